$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 579, shifting the existing
# data block (579:704) down to (581:706).
$ws.Rows("579:580").Insert()

# New row 579 - "Primera" quality entry
$ws.Range("A579").Value = 8
$ws.Range("B579").Value = "Terminal La Palmera de La Serena"
$ws.Range("C579").Value = "Coquimbo"
$ws.Range("D579").Value = 44711
$ws.Range("D579").NumberFormat = $ws.Range("D581").NumberFormat
$ws.Range("E579").Value = 4
$ws.Range("F579").Value = 100112008
$ws.Range("G579").Value = "Coliflor"
$ws.Range("H579").Value = "Sin especificar"
$ws.Range("I579").Value = "Primera"
$ws.Range("J579").Value = 2520
$ws.Range("K579").Value = 800
$ws.Range("L579").Value = 900
$ws.Range("M579").Value = 850
$ws.Range("N579").Value = "$/unidad"
$ws.Range("O579").Value = "Provincia del Elquí"
$ws.Range("P579").Value = 850
$ws.Range("Q579").Value = 1
$ws.Range("R579").Value = "Hortaliza"

# New row 580 - "Segunda" quality entry
$ws.Range("A580").Value = 8
$ws.Range("B580").Value = "Terminal La Palmera de La Serena"
$ws.Range("C580").Value = "Coquimbo"
$ws.Range("D580").Value = 44711
$ws.Range("D580").NumberFormat = $ws.Range("D581").NumberFormat
$ws.Range("E580").Value = 4
$ws.Range("F580").Value = 100112008
$ws.Range("G580").Value = "Coliflor"
$ws.Range("H580").Value = "Sin especificar"
$ws.Range("I580").Value = "Segunda"
$ws.Range("J580").Value = 1400
$ws.Range("K580").Value = 700
$ws.Range("L580").Value = 750
$ws.Range("M580").Value = 725
$ws.Range("N580").Value = "$/unidad"
$ws.Range("O580").Value = "Provincia del Elquí"
$ws.Range("P580").Value = 725
$ws.Range("Q580").Value = 1
$ws.Range("R580").Value = "Hortaliza"
